$d = $word.ActiveDocument

# 1) Title paragraph: drop the leading run(s) of pure whitespace before
#    "Public Health Awareness" (a run of 15 spaces + a single-space
#    "BookTitle" styled run).
$p1 = $d.Paragraphs(1).Range
$titleText = $p1.Text
$idx = $titleText.IndexOf("Public Health")
if ($idx -gt 0) {
    $lead = $d.Range($p1.Start, $p1.Start + $idx)
    $lead.Delete()
}

# 2) Drop the stray trailing space run right after the bold "Conclusion:"
#    heading (it is its own, unformatted run with just " ").
$d.Content.Find.Execute("Conclusion: ", $true, $false, $false, $false, $false, $true, 1, $false, "Conclusion:", 2) | Out-Null

# 3) Swap the author's name from "Rasika M" to "Shobana. M".
$d.Content.Find.Execute("Rasika M", $true, $false, $false, $false, $false, $true, 1, $false, "Shobana. M", 2) | Out-Null
